# Updated excel file to have header column for apps
#
# The sheet originally had column A = app display name (e.g. "Visual C
# Redistributable") and column B = the internal app_var identifier
# (e.g. "visual_c"), with header cells A1 = "app_name" and B1 = "app_var".
# This edit swaps columns A and B (for every data row, including the
# header row) so the internal identifier column becomes column A and the
# human readable name becomes column B.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($lastRow -lt 1) { $lastRow = 32 }

for ($r = 1; $r -le $lastRow; $r++) {
    $colA = $ws.Cells.Item($r, 1)
    $colB = $ws.Cells.Item($r, 2)

    $aVal = $colA.Value2
    $bVal = $colB.Value2

    $colA.Value = $bVal
    $colB.Value = $aVal
}

# Column widths: after the swap, column A holds the short identifiers
# (narrower) and column B holds the longer human readable names (wider),
# so re-fit the two columns to their new best-fit widths.
$ws.Columns.Item(1).ColumnWidth = 17.42
$ws.Columns.Item(2).ColumnWidth = 24.33
